$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before writing, so numeric-looking
# strings like "163.71" are not auto-converted to floating point numbers by
# Excel COM automation (matches original inlineStr/text cell semantics).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.281.39'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '3.136.31'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '571.62'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '163.71'
$ws.Range("E6").Value = '  -2.67%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '0.575'
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("D9").Value = '3.151.59'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("E11").Value = '  -2.81%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '3.688.03'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '64.308.29'
$ws.Range("D16").Value = '25.04'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = '3.145.16'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("E18").Value = '  -2.58%  '
$ws.Range("D19").Value = '403.32'
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("D20").Value = '5.25'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("D21").Value = '12.53'
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '68.95'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("D25").Value = '0.485'
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  -4.81%  '
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '1.80'
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("D32").Value = '21.24'
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("D33").Value = '161.34'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '4.86'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("D35").Value = '6.27'
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("D39").Value = '2.645.82'
$ws.Range("E39").Value = '  -2.48%  '
$ws.Range("D40").Value = '23.70'
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").Value = '4.08'
$ws.Range("E41").Value = '  -2.68%  '
$ws.Range("D42").Value = '38.46'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").Value = '0.689'
$ws.Range("E43").Value = '  -3.53%  '
$ws.Range("D44").Value = '0.0613'
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").Value = '5.41'
$ws.Range("E45").Value = '  -3.15%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0254'
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '21.15'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").Value = '286.86'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '0.0977'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '10.49'
$ws.Range("E51").Value = '  +0.58%  '

# Restore default "Normal" style on column D so no stray style index is
# left attached to the cells (keeps cells styleless, like the original).
$ws.Range("D2:D51").Style = "Normal"

